# Add a new "canonical SMILES" column (D) to the microstates sheet,
# mirroring the "canonical isomeric SMILES" column (C) for each microstate,
# as part of the SAMPL6 microstates update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column D
$ws.Range("D2").Value = "canonical SMILES"

# Populate column D with the same SMILES text already present in column C
# for each of the four microstate rows. (.Value2 avoids the COM property
# wrapper oddities seen with plain .Value when reading back a string.)
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2

# Copy the formatting (fill/font/border/alignment) from column C onto the
# corresponding new column D cells so the new column matches the existing
# banded row styling (header style vs. alternating microstate row styles).
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Give the new column a sensible width (close to the authored 36.86 chars).
$ws.Columns.Item(4).ColumnWidth = 36
